$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '93.648.42'
Set-TextValue 2 5 '  +3.64%  '
Set-TextValue 3 4 '3.133.18'
Set-TextValue 3 5 '  +0.41%  '
Set-TextValue 4 5 '  +0.09%  '
Set-TextValue 5 4 '244.10'
Set-TextValue 5 5 '  +2.89%  '
Set-TextValue 6 4 '617.72'
Set-TextValue 6 5 '  -0.06%  '
Set-TextValue 7 5 '  +1.72%  '
Set-TextValue 8 5 '  +11.57%  '
Set-TextValue 9 4 '0.999'
Set-TextValue 9 5 '  -0.04%  '
Set-TextValue 10 4 '3.130.10'
Set-TextValue 10 5 '  +30.90%  '
Set-TextValue 11 4 '0.749'
Set-TextValue 11 5 '  +1.24%  '
Set-TextValue 12 5 '  -0.29%  '
Set-TextValue 13 5 '  +5.09%  '
Set-TextValue 14 4 '34.94'
Set-TextValue 14 5 '  -0.08%  '
Set-TextValue 15 4 '93.279.72'
Set-TextValue 15 5 '  +3.23%  '
Set-TextValue 16 5 '  -0.06%  '
Set-TextValue 17 4 '3.728.38'
Set-TextValue 18 5 '  -1.73%  '
Set-TextValue 19 5 '  +4.43%  '
Set-TextValue 20 4 '14.95'
Set-TextValue 20 5 '  -0.19%  '
Set-TextValue 21 4 '0.0000210'
Set-TextValue 21 5 '  +4.29%  '
Set-TextValue 22 4 '5.88'
Set-TextValue 22 5 '  +0.69%  '
Set-TextValue 23 2 'Uniswap'
Set-TextValue 23 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 23 4 '9.45'
Set-TextValue 23 5 '  +4.97%  '
Set-TextValue 24 2 'BitcoinCash'
Set-TextValue 24 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 24 4 '451.54'
Set-TextValue 24 5 '  +2.88%  '
Set-TextValue 25 5 '  -1.18%  '
Set-TextValue 26 5 '  +0.70%  '
Set-TextValue 27 4 '11.92'
Set-TextValue 27 5 '  +1.14%  '
Set-TextValue 28 4 '3.302.45'
Set-TextValue 29 5 '  +0.04%  '
Set-TextValue 30 5 '  +10.89%  '
Set-TextValue 31 2 'Stellar'
Set-TextValue 31 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 31 4 '0.230'
Set-TextValue 31 5 '  +0.93%  '
Set-TextValue 32 2 'Cronos'
Set-TextValue 32 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 32 4 '0.171'
Set-TextValue 32 5 '  +1.07%  '
Set-TextValue 33 5 '  +0.65%  '
Set-TextValue 34 5 '  +8.29%  '
Set-TextValue 35 4 '8.14'
Set-TextValue 36 5 '  -4.23%  '
Set-TextValue 37 4 '26.45'
Set-TextValue 37 5 '  +1.51%  '
Set-TextValue 38 5 '  +0.20%  '
Set-TextValue 39 2 'Bittensor'
Set-TextValue 39 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 39 4 '489.15'
Set-TextValue 39 5 '  -3.00%  '
Set-TextValue 40 2 'MantraDAO'
Set-TextValue 40 3 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue 40 4 '3.87'
Set-TextValue 40 5 '  +4.83%  '
Set-TextValue 41 5 '  -2.57%  '
Set-TextValue 42 4 '3.50'
Set-TextValue 42 5 '  +2.12%  '
Set-TextValue 43 5 '  -1.80%  '
Set-TextValue 44 5 '  +4.57%  '
Set-TextValue 45 5 '  +0.00%  '
Set-TextValue 46 4 '163.27'
Set-TextValue 46 5 '  +3.02%  '
Set-TextValue 47 5 '  +2.16%  '
Set-TextValue 48 5 '  -2.74%  '
Set-TextValue 49 4 '1.41'
Set-TextValue 49 5 '  +3.09%  '
Set-TextValue 50 5 '  +5.31%  '
Set-TextValue 51 5 '  +1.06%  '
